# Canada Premier League workbook update (03-06-2024 23:01)
#
# 1. Two already-played fixtures (match ids 7301364 and 6227884, both on
#    45206.75) were re-ordered in the source feed, so rows 83 and 84 swap
#    their entire data payload (everything except the running index in
#    column A, which stays tied to the row position).
# 2. Four upcoming fixtures that had not been played yet (no FT/HT score
#    columns) are removed from the bottom of the table (old rows 116-119,
#    ids 7802881 / 7802946 / 7802947 / 7803370). Removing them also drops
#    those four now-unused id strings from the shared-string table and
#    shrinks the sheet dimension down to AD115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the contents of rows 83 and 84 (columns B through AD) ---
$row83 = $ws.Range("B83:AD83")
$row84 = $ws.Range("B84:AD84")

$vals83 = $row83.Value()
$vals84 = $row84.Value()

$row83.Value = $vals84
$row84.Value = $vals83

# --- 2. Remove the four not-yet-played fixtures at the bottom of the sheet ---
$ws.Range("A116:A119").EntireRow.Delete()
